# Update SYFAFE input cells (rows 6 and 16) with Jun's revised figures.
# All other changed cells in the workbook are formulas that recalc
# automatically from these inputs.
# (Values written in plain decimal form - the PS interpreter here does not
# accept scientific-notation numeric literals like 1.23E-4.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SYFAFE")

# Row 6
$ws.Range("B6").Value = 0.0012523003446111029
$ws.Range("C6").Value = 0.00043095218264560222
$ws.Range("D6").Value = 0.00043095218264560222
$ws.Range("E6").Value = 0.00043095218264560222
$ws.Range("H6").Value = 0.0012928565479368064

# Row 16
$ws.Range("B16").Value = 0.05303186743167481
$ws.Range("C16").Value = 0.016502967462247921
$ws.Range("D16").Value = 0.016502967462247921
$ws.Range("E16").Value = 0.016502967462247921
$ws.Range("H16").Value = 0.049508902386743756
